# Auto-applies the Malboro_Profits.xlsx value updates from the scheduled runner diff.
# For each affected Leve row, refresh the market-price columns (H:N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6: Days of Chunder
$ws.Range("H6").Value = 29.818182
$ws.Range("I6").Value = 29.818182
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 89.45454599999999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 22.54545400000001
$ws.Range("N6").ClearContents()

# ALC row 17: One for the Road
$ws.Range("H17").Value = 2615944.8
$ws.Range("I17").Value = 1997
$ws.Range("J17").Value = 5012063.5
$ws.Range("K17").Value = 5991
$ws.Range("L17").Value = 15036190.5
$ws.Range("M17").Value = -5823
$ws.Range("N17").Value = -15036526.5

# ALC row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 921.8889
$ws.Range("I80").Value = 842.5714
$ws.Range("J80").Value = 1199.5
$ws.Range("K80").Value = 2527.7142
$ws.Range("L80").Value = 3598.5
$ws.Range("M80").Value = -1529.7142
$ws.Range("N80").Value = -5594.5

# ALC row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 921.8889
$ws.Range("I83").Value = 842.5714
$ws.Range("J83").Value = 1199.5
$ws.Range("K83").Value = 7583.1426
$ws.Range("L83").Value = 10795.5
$ws.Range("M83").Value = -2591.1426
$ws.Range("N83").Value = -20779.5

# ALC row 98: The Dotted Line
$ws.Range("H98").Value = 1550.5625
$ws.Range("I98").Value = 1384.1072
$ws.Range("J98").Value = 2715.75
$ws.Range("K98").Value = 1384.1072
$ws.Range("L98").Value = 2715.75
$ws.Range("M98").Value = 113.8928000000001
$ws.Range("N98").Value = -5711.75

# ALC row 106: Making Your Mark
$ws.Range("H106").Value = 12399.111
$ws.Range("I106").Value = 21669
$ws.Range("J106").Value = 811.75
$ws.Range("K106").Value = 21669
$ws.Range("L106").Value = 811.75
$ws.Range("M106").Value = -21038
$ws.Range("N106").Value = -2073.75

# ALC row 122: Wishful Inking
$ws.Range("H122").Value = 1550.5625
$ws.Range("I122").Value = 1384.1072
$ws.Range("J122").Value = 2715.75
$ws.Range("K122").Value = 4152.321599999999
$ws.Range("L122").Value = 8147.25
$ws.Range("M122").Value = -1702.321599999999
$ws.Range("N122").Value = -13047.25

# ALC row 138: All-night Crafting
$ws.Range("H138").Value = 2603.8164
$ws.Range("I138").Value = 2341.1875
$ws.Range("J138").Value = 3098.1765
$ws.Range("K138").Value = 7023.5625
$ws.Range("L138").Value = 9294.529500000001
$ws.Range("M138").Value = -1883.5625
$ws.Range("N138").Value = -19574.5295

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust
$ws.Range("H32").Value = 14039.86
$ws.Range("I32").Value = 4815.415
$ws.Range("J32").Value = 32843.54
$ws.Range("K32").Value = 4815.415
$ws.Range("L32").Value = 32843.54
$ws.Range("M32").Value = -4528.415
$ws.Range("N32").Value = -33417.54

# ARM row 63: Rivets Run through It
$ws.Range("H63").Value = 2941.1667
$ws.Range("I63").Value = 2724.25
$ws.Range("J63").Value = 3375
$ws.Range("K63").Value = 2724.25
$ws.Range("L63").Value = 3375
$ws.Range("M63").Value = -2038.25
$ws.Range("N63").Value = -4747

# ARM row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2941.1667
$ws.Range("I66").Value = 2724.25
$ws.Range("J66").Value = 3375
$ws.Range("K66").Value = 13621.25
$ws.Range("L66").Value = 16875
$ws.Range("M66").Value = -10189.25
$ws.Range("N66").Value = -23739

# ARM row 74: As the Bolt Flies
$ws.Range("H74").Value = 14188.111
$ws.Range("I74").Value = 3377.4443
$ws.Range("J74").Value = 24998.777
$ws.Range("K74").Value = 3377.4443
$ws.Range("L74").Value = 24998.777
$ws.Range("M74").Value = -2503.4443
$ws.Range("N74").Value = -26746.777

# ARM row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 14188.111
$ws.Range("I77").Value = 3377.4443
$ws.Range("J77").Value = 24998.777
$ws.Range("K77").Value = 16887.2215
$ws.Range("L77").Value = 124993.885
$ws.Range("M77").Value = -12519.2215
$ws.Range("N77").Value = -133729.885

# ARM row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 15226.929
$ws.Range("I102").Value = 1013.6923
$ws.Range("J102").Value = 199999
$ws.Range("K102").Value = 1013.6923
$ws.Range("L102").Value = 199999
$ws.Range("M102").Value = 608.3077
$ws.Range("N102").Value = -203243

# ARM row 110: Scheduled Maintenance
$ws.Range("H110").Value = 6896.125
$ws.Range("I110").Value = 8711
$ws.Range("J110").Value = 1451.5
$ws.Range("K110").Value = 8711
$ws.Range("L110").Value = 1451.5
$ws.Range("M110").Value = -6666
$ws.Range("N110").Value = -5541.5

# ARM row 122: Haste for High Durium
$ws.Range("H122").Value = 3464.75
$ws.Range("I122").Value = 2271.2632
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 6813.7896
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -4363.7896
$ws.Range("N122").Value = -28900

# ARM row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2448459.2
$ws.Range("I132").Value = 4404.636
$ws.Range("J132").Value = 5278417
$ws.Range("K132").Value = 13213.908
$ws.Range("L132").Value = 15835251
$ws.Range("M132").Value = -10683.908
$ws.Range("N132").Value = -15840311

$ws = $wb.Worksheets.Item("BSM")
# BSM row 99: Meddle in Metal
$ws.Range("H99").Value = 957.6786
$ws.Range("I99").Value = 882.03705
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 882.03705
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = 615.96295
$ws.Range("N99").Value = -5996

# BSM row 134: Ruthenium Supremium
$ws.Range("H134").Value = 19555.055
$ws.Range("I134").Value = 9905.817999999999
$ws.Range("J134").Value = 34718.145
$ws.Range("K134").Value = 29717.454
$ws.Range("L134").Value = 104154.435
$ws.Range("M134").Value = -27182.454
$ws.Range("N134").Value = -109224.435

$ws = $wb.Worksheets.Item("CRP")
# CRP row 97: Wood That You Could
$ws.Range("H97").Value = 17000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 17000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 17000
$ws.Range("N97").Value = -18982

# CRP row 132: Hull Lotta Damage
$ws.Range("H132").Value = 6425.2583
$ws.Range("I132").Value = 2171.35
$ws.Range("J132").Value = 14159.637
$ws.Range("K132").Value = 6514.049999999999
$ws.Range("L132").Value = 42478.911
$ws.Range("M132").Value = -3984.049999999999
$ws.Range("N132").Value = -47538.911

$ws = $wb.Worksheets.Item("CUL")
# CUL row 86: Let's Not Get Sappy
$ws.Range("H86").Value = 648.62164
$ws.Range("I86").Value = 698.86365
$ws.Range("J86").Value = 574.93335
$ws.Range("K86").Value = 2096.59095
$ws.Range("L86").Value = 1724.80005
$ws.Range("M86").Value = -910.5909499999998
$ws.Range("N86").Value = -4096.80005

# CUL row 89: Luxury Spillover (L)
$ws.Range("H89").Value = 648.62164
$ws.Range("I89").Value = 698.86365
$ws.Range("J89").Value = 574.93335
$ws.Range("K89").Value = 6289.77285
$ws.Range("L89").Value = 5174.40015
$ws.Range("M89").Value = -361.7728500000003
$ws.Range("N89").Value = -17030.40015

# CUL row 133: Friends Are Food
$ws.Range("H133").Value = 4969.778
$ws.Range("I133").Value = 3121.5
$ws.Range("J133").Value = 8666.333000000001
$ws.Range("K133").Value = 9364.5
$ws.Range("L133").Value = 25998.999
$ws.Range("M133").Value = -4304.5
$ws.Range("N133").Value = -36118.999

# CUL row 137: Creative Chocolate
$ws.Range("H137").Value = 2750
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -900
$ws.Range("N137").Value = -20700

# CUL row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 1843.4706
$ws.Range("I140").Value = 1236.2142
$ws.Range("J140").Value = 4677.3335
$ws.Range("K140").Value = 3708.6426
$ws.Range("L140").Value = 14032.0005
$ws.Range("M140").Value = 1471.3574
$ws.Range("N140").Value = -24392.0005

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80: Needs More Prayerbell
$ws.Range("H80").Value = 11128.381
$ws.Range("I80").Value = 9041.166999999999
$ws.Range("J80").Value = 13911.333
$ws.Range("K80").Value = 9041.166999999999
$ws.Range("L80").Value = 13911.333
$ws.Range("M80").Value = -8043.166999999999
$ws.Range("N80").Value = -15907.333

# GSM row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 11128.381
$ws.Range("I83").Value = 9041.166999999999
$ws.Range("J83").Value = 13911.333
$ws.Range("K83").Value = 45205.835
$ws.Range("L83").Value = 69556.66500000001
$ws.Range("M83").Value = -40213.835
$ws.Range("N83").Value = -79540.66500000001

# GSM row 87: Embroiling Embroidery
$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 50000
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

# GSM row 90: The Lovely Hands of Haillenarte (L)
$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 50000
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

# GSM row 132: On Board for Lar
$ws.Range("H132").Value = 5577.394
$ws.Range("I132").Value = 2568.5715
$ws.Range("J132").Value = 10842.833
$ws.Range("K132").Value = 7705.7145
$ws.Range("L132").Value = 32528.499
$ws.Range("M132").Value = -5175.7145
$ws.Range("N132").Value = -37588.499

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40: Best Served Toad
$ws.Range("H40").Value = 6831.6562
$ws.Range("I40").Value = 5415.0625
$ws.Range("J40").Value = 8248.25
$ws.Range("K40").Value = 5415.0625
$ws.Range("L40").Value = 8248.25
$ws.Range("M40").Value = -5279.0625
$ws.Range("N40").Value = -8520.25

# LTW row 82: Trainin' the Neck
$ws.Range("H82").Value = 8264.866
$ws.Range("I82").Value = 8795.799999999999
$ws.Range("J82").Value = 7999.4
$ws.Range("K82").Value = 8795.799999999999
$ws.Range("L82").Value = 7999.4
$ws.Range("M82").Value = -8434.799999999999
$ws.Range("N82").Value = -8721.4

# LTW row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 8264.866
$ws.Range("I85").Value = 8795.799999999999
$ws.Range("J85").Value = 7999.4
$ws.Range("K85").Value = 8795.799999999999
$ws.Range("L85").Value = 7999.4
$ws.Range("M85").Value = -7547.799999999999
$ws.Range("N85").Value = -10495.4

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 1025.3846
$ws.Range("I81").Value = 1025.3846
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2050.7692
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -989.7692000000002

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 1025.3846
$ws.Range("I84").Value = 1025.3846
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 10253.846
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -4949.846000000001

# WVR row 126: A Polished Purchase
$ws.Range("H126").Value = 3819.652
$ws.Range("I126").Value = 3881.5
$ws.Range("J126").Value = 3407.3333
$ws.Range("K126").Value = 11644.5
$ws.Range("L126").Value = 10221.9999
$ws.Range("M126").Value = -9174.5
$ws.Range("N126").Value = -15161.9999

# WVR row 132: Comfy Cabins
$ws.Range("H132").Value = 7089.1562
$ws.Range("I132").Value = 2283.8096
$ws.Range("J132").Value = 16263
$ws.Range("K132").Value = 6851.4288
$ws.Range("L132").Value = 48789
$ws.Range("M132").Value = -4321.4288
$ws.Range("N132").Value = -53849
